$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and the two row-pair
# reorderings for rows 48-51) as captured by the source diff.
# NumberFormat is forced to Text ("@") before writing so that purely
# numeric-looking values (e.g. "237.36") are stored as inline strings,
# matching the original sheet layout, rather than being auto-converted
# to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "94.371.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.427.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.36"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -5.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "641.89"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.53%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.405"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.04%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.967"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -7.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.426.55"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.57"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.18"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.187.40"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.068.95"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000251"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.29"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -6.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.425.74"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.47"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.49"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.494"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -8.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "497.46"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.23"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.85%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.48"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -8.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "93.73"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.97"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.609.19"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.65"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.76"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.138"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "29.70"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.552"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "544.40"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.64"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.45"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.906"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.06"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.34"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.63%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "MantraDAO"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.57"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0408"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.55%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.20"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.54%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.88"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.06%  "
